# Apply the "walkingToRunning" sensor-data update.
#
# The new capture window slides the data forward: what used to be rows
# 11-21 (timestamps 900-1900) becomes the new rows 2-12, and 10 fresh
# samples (timestamps 2000-2900) are appended as rows 22-31. The A
# (timestamp) and B (label) columns are untouched by this refresh - only
# the sensor columns C:H move/change, plus the new rows 22-31 are added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement block for A2:H31 (30 rows x 8 cols) - timestamp, label,
# and the 6 refreshed sensor columns (ax, ay, az, gx, gy, gz).
$rows = 30
$cols = 8
$data = New-Object 'object[,]' $rows,$cols

# row 2  (timestamp 0)
$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = 2.993003823075968
$data[0,3] = 3.121233616556498
$data[0,4] = -0.8490933929171138
$data[0,5] = 2.092578887939453
$data[0,6] = -8.521716117858887
$data[0,7] = 1.690641283988953
# row 3  (timestamp 100)
$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = -4.831494113377142
$data[1,3] = -9.479980850219709
$data[1,4] = 16.28942495073587
$data[1,5] = 2.046772480010986
$data[1,6] = -4.874783515930176
$data[1,7] = 0.8943560719490051
# row 4  (timestamp 200)
$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = -9.035116481781008
$data[2,3] = -3.15880911350261
$data[2,4] = 10.23342611789717
$data[2,5] = 4.356798648834229
$data[2,6] = -5.594369411468506
$data[2,7] = -4.091081619262695
# row 5  (timestamp 300)
$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = 7.237885883876157
$data[3,3] = -10.53579415593821
$data[3,4] = 15.53679728507988
$data[3,5] = -5.464406967163086
$data[3,6] = 3.147723197937012
$data[3,7] = 2.375073671340942
# row 6  (timestamp 400)
$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = -46.59491330555468
$data[4,3] = -20.57918526785712
$data[4,4] = -11.76733117784752
$data[4,5] = -6.766693115234375
$data[4,6] = 7.695337295532227
$data[4,7] = 7.492071628570557
# row 7  (timestamp 500)
$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = -38.6563062940324
$data[5,3] = -8.040147893769365
$data[5,4] = -22.76833855765205
$data[5,5] = 3.930692672729492
$data[5,6] = 9.625063896179199
$data[5,7] = -5.505752563476562
# row 8  (timestamp 600)
$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = -17.68361431530541
$data[6,3] = -12.95902962642049
$data[6,4] = -10.52790222508571
$data[6,5] = 3.03800106048584
$data[6,6] = -1.461676001548767
$data[6,7] = 0.6908905506134033
# row 9  (timestamp 700)
$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = -13.64767029626039
$data[7,3] = -8.337594250270474
$data[7,4] = 4.581750222614684
$data[7,5] = 11.57449817657471
$data[7,6] = -5.450558662414551
$data[7,7] = -0.0244345031678676
# row 10  (timestamp 800)
$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = 12.93104904038545
$data[8,3] = 0.6822976725441974
$data[8,4] = 10.89643403462002
$data[8,5] = -4.411925792694092
$data[8,6] = 2.68752908706665
$data[8,7] = 3.053647041320801
# row 11  (timestamp 900)
$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = -0.9200370788569145
$data[9,3] = -13.92125396728488
$data[9,4] = -0.9832572937009605
$data[9,5] = -4.802345275878906
$data[9,6] = 4.841493606567383
$data[9,7] = 2.432597875595093
# row 12  (timestamp 1000)
$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = -8.313682734966347
$data[10,3] = -9.559037576615895
$data[10,4] = -17.89249617712833
$data[10,5] = -2.648380517959595
$data[10,6] = 16.29841232299805
$data[10,7] = -2.319014072418213
# row 13  (timestamp 1100)
$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = -8.376615575381617
$data[11,3] = -8.225071150490136
$data[11,4] = -15.93280724116742
$data[11,5] = -0.7825698852539062
$data[11,6] = -1.549027681350708
$data[11,7] = 1.645367503166199
# row 14  (timestamp 1200)
$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = 0.8738524981906739
$data[12,3] = -4.53215429782872
$data[12,4] = 15.28993327277035
$data[12,5] = 9.456219673156738
$data[12,6] = -6.493985176086426
$data[12,7] = -6.240252494812012
# row 15  (timestamp 1300)
$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = -0.01965822492311098
$data[13,3] = 16.21388537543142
$data[13,4] = 17.33802195957751
$data[13,5] = 3.910985231399536
$data[13,6] = 0.8456867933273315
$data[13,7] = -5.051949977874756
# row 16  (timestamp 1400)
$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = -10.3103993143354
$data[14,3] = -31.18352238791278
$data[14,4] = 18.67163455826883
$data[14,5] = -2.794321775436401
$data[14,6] = -2.126400947570801
$data[14,7] = -1.946171641349792
# row 17  (timestamp 1500)
$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = 22.02176739828893
$data[15,3] = -22.72049116407177
$data[15,4] = -2.185512433732939
$data[15,5] = -14.5798749923706
$data[15,6] = -11.00378227233887
$data[15,7] = 6.221210956573486
# row 18  (timestamp 1600)
$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = 5.463320446014379
$data[16,3] = 6.157036137580864
$data[16,4] = -27.94153518676758
$data[16,5] = 5.833254814147949
$data[16,6] = -11.42988777160644
$data[16,7] = 6.618554592132568
# row 19  (timestamp 1700)
$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = -37.06017020089279
$data[17,3] = -23.85922709873728
$data[17,4] = 0.4766027178080634
$data[17,5] = 0.6635265946388245
$data[17,6] = -4.405001640319824
$data[17,7] = 1.273590207099915
# row 20  (timestamp 1800)
$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = -19.13142990640274
$data[18,3] = -10.62035506112288
$data[18,4] = 8.715779474803668
$data[18,5] = 5.506218910217285
$data[18,6] = 16.36765480041504
$data[18,7] = -1.815144062042236
# row 21  (timestamp 1900)
$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = 6.810042284216134
$data[19,3] = 8.792214941978534
$data[19,4] = -5.953527409689759
$data[19,5] = -1.695501565933228
$data[19,6] = -2.446512937545776
$data[19,7] = 2.133258581161499
# row 22  (timestamp 2000)
$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = -0.03172696658525354
$data[20,3] = -6.182636371680525
$data[20,4] = -11.13081671851015
$data[20,5] = -4.14081621170044
$data[20,6] = 5.169595241546631
$data[20,7] = 1.407280921936035
# row 23  (timestamp 2100)
$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = 12.65357358115054
$data[21,3] = -1.313164901733439
$data[21,4] = -15.5251411982945
$data[21,5] = -2.688860654830933
$data[21,6] = 10.50284194946289
$data[21,7] = -1.326720356941223
# row 24  (timestamp 2200)
$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = 8.714159175328177
$data[22,3] = -9.301818609237403
$data[22,4] = -17.37430092947822
$data[22,5] = -1.443033814430237
$data[22,6] = 10.46609020233154
$data[22,7] = -6.833072185516357
# row 25  (timestamp 2300)
$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = -5.424608945846559
$data[23,3] = -50.41501617431641
$data[23,4] = -2.802897453308105
$data[23,5] = 1.430516958236694
$data[23,6] = -2.983938932418823
$data[23,7] = -5.8599534034729
# row 26  (timestamp 2400)
$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = 3.029158597333116
$data[24,3] = -28.14689919608023
$data[24,4] = 5.423780044487627
$data[24,5] = 3.719237804412842
$data[24,6] = 1.893374443054199
$data[24,7] = -7.143064022064209
# row 27  (timestamp 2500)
$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = 23.23876049859135
$data[25,3] = 5.670972211020553
$data[25,4] = 3.751934364863664
$data[25,5] = -0.8805742263793945
$data[25,6] = -0.7974836230278015
$data[25,7] = 1.437640905380249
# row 28  (timestamp 2600)
$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = 12.21474557604104
$data[26,3] = -7.668160905156817
$data[26,4] = -2.514941801343662
$data[26,5] = -4.217514991760254
$data[26,6] = 1.791641712188721
$data[26,7] = 6.048105716705322
# row 29  (timestamp 2700)
$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = 1.01306713649177
$data[27,3] = 2.441266550336438
$data[27,4] = -15.73728098188129
$data[27,5] = -1.25075364112854
$data[27,6] = 12.15346908569336
$data[27,7] = -1.156278014183044
# row 30  (timestamp 2800)
$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -0.6987579890658431
$data[28,3] = -8.856165681566392
$data[28,4] = -19.26768711635043
$data[28,5] = -1.032907009124756
$data[28,6] = 9.82000732421875
$data[28,7] = -6.767558574676514
# row 31  (timestamp 2900)
$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = -16.89904033115963
$data[29,3] = -34.02156884329665
$data[29,4] = 7.39472787039573
$data[29,5] = -3.240667581558228
$data[29,6] = -12.29328536987305
$data[29,7] = -4.838364601135254

$ws.Range("A2:H31").Value = $data

